$d = $word.ActiveDocument

# Update the date line (first paragraph)
$d.Paragraphs.Item(1).Range.Text = "2024-05-19 Sunday"

# Update each arithmetic cell in the 20x5 table, in row-major order
$tbl = $d.Tables.Item(1)
$values = @(
    "74-60=",
    "66-29=",
    "94-41=",
    "75+16=",
    "84-51=",
    "19+32=",
    "4+87=",
    "32+32=",
    "30-5=",
    "58+39=",
    "57+39=",
    "2+23=",
    "31-27=",
    "17+28=",
    "4+37=",
    "87-1=",
    "53-8=",
    "81+3=",
    "37-0=",
    "8+67=",
    "53+8=",
    "63-61=",
    "42-23=",
    "57+26=",
    "15+25=",
    "7+47=",
    "26+48=",
    "89-43=",
    "34+40=",
    "39+37=",
    "9+87=",
    "77-8=",
    "2+50=",
    "16+44=",
    "32+23=",
    "79-8=",
    "79-18=",
    "89-3=",
    "74+13=",
    "2+78=",
    "62+27=",
    "13+52=",
    "49-6=",
    "46+20=",
    "92+4=",
    "65+10=",
    "83-8=",
    "69-67=",
    "8+10=",
    "83-35=",
    "47+52=",
    "90-13=",
    "35+28=",
    "88-63=",
    "88-77=",
    "4+16=",
    "44+34=",
    "80-35=",
    "49-5=",
    "68-5=",
    "48+4=",
    "77-7=",
    "75-48=",
    "6+26=",
    "49-7=",
    "46+43=",
    "36+50=",
    "68+4=",
    "66-47=",
    "40-0=",
    "92-84=",
    "42-19=",
    "0+45=",
    "20+60=",
    "62-10=",
    "66+27=",
    "51+5=",
    "38-27=",
    "24+34=",
    "59-45=",
    "31-18=",
    "58+40=",
    "78-7=",
    "35-30=",
    "57+41=",
    "13+31=",
    "39+24=",
    "4+72=",
    "82-78=",
    "4+42=",
    "78+0=",
    "95-43=",
    "63+28=",
    "9+16=",
    "71+9=",
    "3+7=",
    "37+38=",
    "37+48=",
    "72-28=",
    "92-88="
)

$nCols = 5
$idx = 0
for ($r = 1; $r -le 20; $r++) {
    for ($c = 1; $c -le $nCols; $c++) {
        $cell = $tbl.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
        $idx++
    }
}

Write-Output "done"
